$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value2 = $ws.Range("A1").Value2
$ws.Range("D2").Value2 = $ws.Range("A2").Value2.Replace("-", "")
$ws.Range("D3").Value2 = $ws.Range("A3").Value2.Replace("-", "")
$ws.Range("D4").Value2 = $ws.Range("A4").Value2.Replace("-", "")

$ws.Range("D4").Select()
